# Update the date heading
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-11-02 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-11-03 Sunday", 2)

# Update the division problems in the table. The table has empty spacer
# rows interleaved with the 5 data rows, so addressing cells directly by
# (row, column) avoids any ambiguity from duplicate problem text (e.g.
# "11÷5=" and "74÷2=" each appear twice but map to different results).
$t = $d.Tables.Item(1)

$values = @(
    @(1,  1, "88÷7="),
    @(1,  2, "54÷7="),
    @(1,  3, "36÷8="),
    @(1,  4, "37÷6="),
    @(1,  5, "60÷2="),

    @(5,  1, "99÷5="),
    @(5,  2, "29÷9="),
    @(5,  3, "43÷3="),
    @(5,  4, "11÷8="),
    @(5,  5, "32÷4="),

    @(9,  1, "45÷4="),
    @(9,  2, "59÷5="),
    @(9,  3, "36÷3="),
    @(9,  4, "55÷2="),
    @(9,  5, "16÷3="),

    @(13, 1, "94÷4="),
    @(13, 2, "96÷7="),
    @(13, 3, "83÷3="),
    @(13, 4, "91÷8="),
    @(13, 5, "83÷4="),

    @(17, 1, "27÷3="),
    @(17, 2, "97÷5="),
    @(17, 3, "76÷8="),
    @(17, 4, "22÷7="),
    @(17, 5, "26÷8=")
)

foreach ($entry in $values) {
    $row = $entry[0]
    $col = $entry[1]
    $newText = $entry[2]
    $cell = $t.Cell($row, $col)
    $cell.Range.Text = $newText
}
